# uOmodri_v2_BOM.xlsx fixes:
# - Drills too close to the edge board.
# - Mechanical pads of power transistor too close to edge board.
# - Some JTAG tracks didn't respect edge clearance.
#
# This translates into BOM bookkeeping updates: several component rows get a
# STOCK (column J) quantity filled in, a couple of STOCK values are updated
# to reflect extra stock on hand, and those rows are marked with the
# "Satisfaisant" (green) cell style to show they have been checked/resolved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 already uses the plain "Satisfaisant" style (green fill, thin border,
# centered/wrapped text) that rows 13, 18, 20, 22, 30, 31 and 38 need to move
# to. Copying its formatting (instead of assigning a named style) makes Excel
# reuse the existing cell style definition rather than create a new one.
$formatDonor = $ws.Range("A14:J14")

# --- Row 13: C59, C98 22uF capacitor now has stock -----------------------
$ws.Range("J13").Value = 50
$formatDonor.Copy()
$ws.Range("A13:J13").PasteSpecial(-4122)

# --- Row 18: 909 resistor now has stock -----------------------------------
$ws.Range("J18").Value = 100
$formatDonor.Copy()
$ws.Range("A18:J18").PasteSpecial(-4122)

# --- Row 20: resistor now has stock ---------------------------------------
$ws.Range("J20").Value = 100
$formatDonor.Copy()
$ws.Range("A20:J20").PasteSpecial(-4122)

# --- Row 22: STOCK raised from "200+" to "400+" ---------------------------
$ws.Range("J22").Value = "400+"
$formatDonor.Copy()
$ws.Range("A22:J22").PasteSpecial(-4122)

# --- Row 30: now has stock --------------------------------------------------
$ws.Range("J30").Value = 25
$formatDonor.Copy()
$ws.Range("A30:J30").PasteSpecial(-4122)

# --- Row 31: STOCK updated from "5" to "5+10" -----------------------------
$ws.Range("J31").Value = "5+10"
$ws.Range("A14:G14").Copy()
$ws.Range("A31:G31").PasteSpecial(-4122)
$ws.Range("I14:J14").Copy()
$ws.Range("I31:J31").PasteSpecial(-4122)
# H31 has no border (unlike the rest of the row); give it the same green
# "Satisfaisant" fill/font/alignment as its neighbours but keep it borderless.
$ws.Range("A14").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Borders.LineStyle = -4142

# --- Row 38: STOCK updated from "11" to "11+10" ---------------------------
$ws.Range("J38").Value = "11+10"
$formatDonor.Copy()
$ws.Range("A38:J38").PasteSpecial(-4122)

# The "Neutre" (orange, warning) cell style is no longer used by any cell now
# that rows 22, 31 and 38 have been marked "Satisfaisant" (green) -> drop it.
$wb.Styles("Neutre").Delete()

# Restore the selection to where the author last left it.
$null = $ws.Range("G19").Select()

$wb.Save()
